$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 81
$ws1.Range("F5").Value = 1716
$ws1.Range("F7").Value = 937
$ws1.Range("F8").Value = 2120
$ws1.Range("F9").Value = 2045
$ws1.Range("F10").Value = 1059
$ws1.Range("F13").Value = 1638
$ws1.Range("F14").Value = 360
$ws1.Range("F16").Value = 20
$ws1.Range("F18").Value = 130
$ws1.Range("F19").Value = 1499
$ws1.Range("F20").Value = 566
$ws1.Range("F21").Value = 666
$ws1.Range("F22").Value = 552
$ws1.Range("F23").Value = 11941
$ws1.Range("F24").Value = 11956
$ws1.Range("F27").Value = 269
$ws1.Range("F29").Value = 172
$ws1.Range("F30").Value = 498

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value = 37
$ws2.Range("F7").Value = 6

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 81
$ws4.Range("F7").Value = 1716
$ws4.Range("F9").Value = 937
$ws4.Range("F10").Value = 2120
$ws4.Range("F11").Value = 2045
$ws4.Range("F12").Value = 1060
$ws4.Range("F15").Value = 1638
$ws4.Range("F16").Value = 360
$ws4.Range("F18").Value = 20
$ws4.Range("F22").Value = 130
$ws4.Range("F23").Value = 1499
$ws4.Range("F24").Value = 566
$ws4.Range("F25").Value = 666
$ws4.Range("F26").Value = 552
$ws4.Range("F27").Value = 11941
$ws4.Range("F28").Value = 11956
$ws4.Range("F31").Value = 269
$ws4.Range("F34").Value = 37
$ws4.Range("F35").Value = 172
$ws4.Range("F36").Value = 498
$ws4.Range("F37").Value = 6
